# "Added training review result for Ruchi ma'am"
#
# Table 3 (rows 21-28) is the "RUCHI PAREEK" review block, scored by
# Priyanka (col C) and Aashna (col D), out of Max. Marks (col E).
# Priyanka's scores (col C) were missing; fill them in, along with the
# "NA" entries + max-marks that accompany them on the last two rows, and
# fix up the column-C total formula to match column D's.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - used to copy an existing cell's style (cellXf) onto a
# target cell without disturbing its own value/formula.
$xlPasteFormats = -4122

# --- Row 23: React Js Understanding / Basic concept understanding ---
# C23 was blank with style 9 ("NA"-style fill); give it the numeric style
# (style 6, same as D23) and the score.
$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial($xlPasteFormats)
$ws.Range("C23").Value = 1.5

# --- Row 24: Presentation / Fluency, Presentation of content ---
$ws.Range("D23").Copy()
$ws.Range("C24").PasteSpecial($xlPasteFormats)
$ws.Range("C24").Value = 1

# --- Row 25: Interview Performance / Understanding, Logical Understanding ---
$ws.Range("D23").Copy()
$ws.Range("C25").PasteSpecial($xlPasteFormats)
$ws.Range("C25").Value = 1.5

# --- Row 26: Project Implementation / Project Setup, Routing, Components ---
# Both C26 and D26 become "NA" (text), with the "NA" fill style (style 9,
# same style C26 already had); E26 gets the max marks value.
$ws.Range("C26").Value = "NA"

$ws.Range("C26").Copy()
$ws.Range("D26").PasteSpecial($xlPasteFormats)
$ws.Range("D26").Value = "NA"

$ws.Range("E26").Value = 2

# --- Row 27: Code Quality / Indentation, Code Quality, Project Setup ---
$ws.Range("C27").Value = "NA"

$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$ws.Range("D27").Value = "NA"

$ws.Range("E27").Value = 2

# --- Row 28: totals ---
# C28 gets the same "sum" style as D28 (and the matching tables' totals)
# plus the mirroring SUM formula.
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial($xlPasteFormats)
$ws.Range("C28").Formula = "=SUM(C23:C27)"

# The "RUCHI PAREEK" section header row grew a bit to fit the now-complete
# table beneath it.
$ws.Rows.Item(21).RowHeight = 24

# --- view state: scrolled down a couple rows, selection moved to G25 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("G25").Select()
